# Apply updated crypto price/volume figures (and one coin swap: Stacks -> EOS)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '28.209.45'
$ws.Range('E2').Value = '  +4.00%  '
# Row 3
$ws.Range('D3').Value = '1.785.24'
$ws.Range('E3').Value = '  +0.12%  '
# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9985'
$ws.Range('E4').Value = '  -0.56%  '
# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '336.46'
$ws.Range('E5').Value = '  -0.01%  '
# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9956'
$ws.Range('E6').Value = '  -0.48%  '
# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3831'
$ws.Range('E7').Value = '  +0.31%  '
# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3438'
$ws.Range('E8').Value = '  +0.86%  '
# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '47.69'
$ws.Range('E9').Value = '  -0.65%  '
# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.161'
$ws.Range('E10').Value = '  -2.22%  '
# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07438'
$ws.Range('E11').Value = '  -0.02%  '
# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '23.20'
$ws.Range('E12').Value = '  +7.01%  '
# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.9957'
$ws.Range('E13').Value = '  -0.64%  '
# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.422'
$ws.Range('E14').Value = '  -0.14%  '
# Row 15
$ws.Range('D15').Value = '1.783.87'
$ws.Range('E15').Value = '  +0.04%  '
# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.134'
$ws.Range('E16').Value = '  +0.65%  '
# Row 17
$ws.Range('E17').Value = '  -0.74%  '
# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06654'
$ws.Range('E18').Value = '  +0.11%  '
# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '82.84'
$ws.Range('E19').Value = '  -0.64%  '
# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.9959'
$ws.Range('E20').Value = '  -0.49%  '
# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.51'
# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.441'
$ws.Range('E22').Value = '  -1.19%  '
# Row 23
$ws.Range('D23').Value = '28.209.81'
$ws.Range('E23').Value = '  +3.95%  '
# Row 24
$ws.Range('E24').Value = '  -1.12%  '
# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.387'
$ws.Range('E25').Value = '  +0.54%  '
# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.444'
$ws.Range('E26').Value = '  +0.00%  '
# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.90'
$ws.Range('E27').Value = '  -0.95%  '
# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.431'
$ws.Range('E28').Value = '  -2.37%  '
# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '153.88'
$ws.Range('E29').Value = '  -0.97%  '
# Row 30
$ws.Range('D30').Value = '1.986.66'
$ws.Range('E30').Value = '  +0.10%  '
# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '134.69'
$ws.Range('E31').Value = '  +0.46%  '
# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.162'
$ws.Range('E32').Value = '  +2.03%  '
# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.962'
$ws.Range('E33').Value = '  -0.66%  '
# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08803'
$ws.Range('E34').Value = '  +1.58%  '
# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '12.79'
$ws.Range('E35').Value = '  -1.82%  '
# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02426'
$ws.Range('E36').Value = '  +4.70%  '
# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.6865'
$ws.Range('E37').Value = '  +0.58%  '
# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.334'
$ws.Range('E38').Value = '  -0.86%  '
# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06341'
$ws.Range('E39').Value = '  +0.85%  '
# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.2185'
$ws.Range('E40').Value = '  +0.41%  '
# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.516'
$ws.Range('E41').Value = '  -6.58%  '
# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.244'
$ws.Range('E42').Value = '  +0.31%  '
# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.343'
$ws.Range('E43').Value = '  +0.28%  '
# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '14.27'
$ws.Range('E44').Value = '  +0.25%  '
# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.9950'
# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6328'
$ws.Range('E46').Value = '  -1.46%  '
# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.851'
$ws.Range('E47').Value = '  -0.07%  '
# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '132.34'
$ws.Range('E48').Value = '  +0.79%  '
# Row 49
$ws.Range('E49').Value = '  -1.44%  '
# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.07445'
$ws.Range('E50').Value = '  +4.81%  '
# Row 51
$ws.Range('B51').Value = 'EOS'
$ws.Range('C51').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.185'
$ws.Range('E51').Value = '  +6.56%  '
